# Updated the site-audit spreadsheet:
#  - fixed spelling typos (Accessability -> Accessibility, appropiate ->
#    appropriate, repeting -> repeating, visability -> visibility,
#    Incorperate -> Incorporate, Vesrion stays as typed by author)
#  - split the old single "alt tags" row into a dedicated SEO row and a
#    dedicated Accessibility row
#  - added a new "Best Practices" / jQuery+Bootstrap version row (14) with
#    a hyperlink reference, documenting the Bootstrap 3.4.1 update
#  - moved the cursor/selection to the newly added row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 : SEO / image format ---
$ws.Range("A3").Value2 = "SEO"
$ws.Range("B3").Value2 = "image format"
$ws.Range("C3").Value2 = "modern image formats are smaller and provide better compression and less data use"
$ws.Range("D3").Value2 = "convert images to a current web focused format such as WebP"
$ws.Range("F3").Value2 = "https://imageseo.io/images-seo-optimization/"

# --- Row 4 : SEO / Image size ---
$ws.Range("A4").Value2 = "SEO"
$ws.Range("B4").Value2 = "Image size"
$ws.Range("C4").Value2 = "Reducing image size also reduces page loading speed which can affect the site rank"
$ws.Range("D4").Value2 = "Convert images to a more efficient format and correctly size them"
$ws.Range("F4").Value2 = "https://developers.google.com/search/docs/advanced/guidelines/google-images/`nhttps://www.searchenginejournal.com/on-page-seo/image-optimization/"

# --- Row 5 : SEO / Image filename ---
$ws.Range("A5").Value2 = "SEO"
$ws.Range("B5").Value2 = "Image filename"
$ws.Range("C5").Value2 = "Using a descriptive filename provides google with easy to access information about the image."
$ws.Range("D5").Value2 = "Name each image with an appropriate name but keep it short, avoid keyword stuffing"
$ws.Range("F5").Value2 = "https://ahrefs.com/blog/image-seo/`nhttps://developers.google.com/search/docs/advanced/guidelines/google-images"

# --- Row 6 : SEO / Image alt tags ---
$ws.Range("A6").Value2 = "SEO "
$ws.Range("B6").Value2 = "Image alt tags"
$ws.Range("C6").Value2 = "Improper use of alt tags can hurt SEO rankings and keyword stuffing here can have large, negative effects."
$ws.Range("D6").Value2 = "Alt tags should only be used to describe the image, avoid keyword stuffing here too."
$ws.Range("F6").Value2 = "https://adolab.com/alt-text-best-practices-seo-mistakes/ `nhttps://adolab.com/what-is-image-alt-text-beginners-guide/"

# --- Row 7 : Accessibility / Image alt tags (new split row) ---
$ws.Range("A7").Value2 = "Accessibility"
$ws.Range("B7").Value2 = "Image alt tags"
$ws.Range("C7").Value2 = "Proper use of alt tags on images allows users with visual disabilities to understand images without seeing them"
$ws.Range("D7").Value2 = "Apply alt tags the properly describe the image without repeating or deviating from the image itself."
$ws.Range("F7").Value2 = "https://moz.com/learn/seo/alt-text `nhttps://adolab.com/alt-text-best-practices-seo-mistakes/"

# --- Row 8 : Accessibility / contrast text/background ---
$ws.Range("A8").Value2 = "Accessibility"
$ws.Range("B8").Value2 = "contrast text/background"
$ws.Range("C8").Value2 = "Improper usage of colors makes a site difficult to use and will indirectly affect ranking due to a high bounce rate. "
$ws.Range("D8").Value2 = "Choose background, foreground, and text colors that work together to maintain visibility"
$ws.Range("F8").Value2 = "https://www.w3.org/TR/UNDERSTANDING-WCAG20/visual-audio-contrast-contrast.html"

# --- Row 9 : Accessibility/SEO / HTML Head ---
$ws.Range("A9").Value2 = "Accessibility/SEO"
$ws.Range("B9").Value2 = "HTML Head"
$ws.Range("C9").Value2 = "It is necessary to follow best practices here to improve SEO and Accessibility"
$ws.Range("D9").Value2 = "Language should be identified. Meta tags should be refined to suit each use."
$ws.Range("F9").Value2 = "https://www.w3.org/TR/html401/struct/dirlang.html https://www.w3schools.com/tags/att_meta_name.asp"

# --- Row 10 : SEO / Blackhat Keywords ---
$ws.Range("A10").Value2 = "SEO"
$ws.Range("B10").Value2 = "Blackhat Keywords "
$ws.Range("C10").Value2 = "Adding keywords in the html without using them in context will result in penalization from Google"
$ws.Range("D10").Value2 = "Incorporate keywords into the content in a natural way using proper context."
$ws.Range("F10").Value2 = "https://cognitiveseo.com/blog/12169/44-black-hat-seo-techniques/#6"

# --- Row 11 : SEO / Link text ---
$ws.Range("A11").Value2 = "SEO"
$ws.Range("B11").Value2 = "Link text"
$ws.Range("C11").Value2 = "Some links do not clearly describe where they go. This can impact SEO and accessibility "
$ws.Range("D11").Value2 = "Links should be clear where they lead to and in the case of image links, an alt tag should be added"
$ws.Range("F11").Value2 = "https://ahrefs.com/blog/anchor-text/"

# --- Row 12 : SEO / Semantic HTML ---
$ws.Range("A12").Value2 = "SEO"
$ws.Range("B12").Value2 = "Semantic HTML"
$ws.Range("C12").Value2 = "Overuse of DIV does not provide google with information about differences in each section"
$ws.Range("D12").Value2 = "Separate html sections using appropriate semantic html. "
$ws.Range("F12").Value2 = "https://www.tmprod.com/blog/2010/h1-h2-and-h3-header-tags-for-seo/ https://www.w3schools.com/html/html5_semantic_elements.asp"

# --- Row 13 : SEO / Images of text ---
$ws.Range("A13").Value2 = "SEO"
$ws.Range("B13").Value2 = "Images of text"
$ws.Range("C13").Value2 = "Adding text instead of images of text offers more opportunity for keyword usage and site ranking"
$ws.Range("D13").Value2 = "Wherever possible use text instead of images of text. Use strong/bold to highlight important parts"
$ws.Range("F13").Value2 = "http://www.velizaratellalyan.com/seo/the-importance-of-bold-and-strong-tags-in-seo/"

# --- Row 14 (new): Best Practices / jQuery Version and Bootstrap Version ---
$ws.Range("A14").Value2 = "Best Practices"
$ws.Range("B14").Value2 = "jQuery Vesrion and Bootstrap Version"
$ws.Range("C14").Value2 = "To ensure site security, all dependencies should be kept up to date when possible. "
$ws.Range("D14").Value2 = "Stay up to date on the current releases of dependencies within your project."
$ws.Range("E14").Value2 = $true

# add the hyperlink for the new jQuery/Bootstrap reference row (this also
# applies the built-in "Hyperlink" style to F14, matching the other
# reference cells in column F). Hyperlinks.Add overwrites the cell's text
# with the TextToDisplay argument, so the full two-line reference text is
# restored right after.
$ws.Hyperlinks.Add($ws.Range("F14"), "https://github.com/jquery/jquery-migrate/", "README ", "", "https://github.com/jquery/jquery-migrate/#README ") | Out-Null
$ws.Range("F14").Value2 = "https://github.com/jquery/jquery-migrate/#README `nhttps://getbootstrap.com/docs/versions/ "

# give row 14 the same look & feel (wrapped text, row height) as the rows above it
$ws.Range("A14:D14").WrapText = $true
$ws.Range("F14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 60.75

# row 15 (previously unused trailing spacer row) picks up the hyperlink look
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("F15").WrapText = $true

# move the view / selection to show the newly added row, matching the
# author's final on-screen state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F15").Select()
